$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook-level: refresh the revisionPtr documentId and the window geometry
# (these are view-only bookkeeping fields Excel rewrites on save/open).
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Left = -120
$wb.Windows.Item(1).Top = -120
$wb.Windows.Item(1).Width = 51840
$wb.Windows.Item(1).Height = 21240

# ---------------------------------------------------------------------------
# SPRINT 1 sheet: add two new work-log entries and re-date / re-size a
# couple of the existing ones.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SPRINT 1")

# Make room: one new row right before the "JTI synthesis" row, and a second
# new row right before the (soon to be overwritten) "conclusion" row.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(10).Insert()
# The blank filler block below the table only needs to grow by one row
# (not two), so trim one row back out of it.
$ws.Rows.Item(16).Delete()

# New row 7: reading done before finishing the JTI synthesis
$ws.Range("B7").Value = 43911
$ws.Range("C7").Value = "Lecture du CWA 17493 JTI"
$ws.Range("D7").Value = 2

# Row 8 (was row 7): JTI synthesis, same text, date pushed to 2020-03-26
$ws.Range("B8").Value = 43916

# Row 9 (was row 8): comparison synthesis, date pushed to 2020-03-26, time updated
$ws.Range("B9").Value = 43916
$ws.Range("D9").Value = 3

# New row 10: page layout + conclusion work
$ws.Range("B10").Value = 43916
$ws.Range("C10").Value = "Mise en page état de l'art et conclusion"
$ws.Range("D10").Value = 1

# Row 11 (was row 9, the old "conclusion" entry): replaced by the Teams
# meeting / sprint review entry
$ws.Range("B11").Value = 43916
$ws.Range("C11").Value = "Rencontre sur Teams + Sprint review et spring 2 planning"
$ws.Range("D11").Value = 1

[void]$ws.Range("D12").Select()
